{"js": "// The document ends with an empty paragraph right before the section break.\n// We need to turn that empty paragraph into \"Student Name: Ekamjot Singh\"\n// and add a new paragraph after it with \"GitHub Username: jotekam4\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertText(\"Student Name: Ekamjot Singh\", Word.InsertLocation.replace);\nlastParagraph.insertParagraph(\"GitHub Username: jotekam4\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# The document ends with an empty paragraph right before the section break.\n# Turn that empty paragraph into \"Student Name: Ekamjot Singh\" and add a new\n# paragraph after it with \"GitHub Username: jotekam4\".\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs($d.Paragraphs.Count)\n$lastParagraph.Range.Text = \"Student Name: Ekamjot Singh\"\n\n$lastParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs($d.Paragraphs.Count)\n$newParagraph.Range.Text = \"GitHub Username: jotekam4\"\n"}
